$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header labels with units
$ws.Range("J1").Value = "MAE [`$COP/kWh]"
$ws.Range("K1").Value = "MSE [`$COP/kWh]"
$ws.Range("L1").Value = "RMSE [`$COP/kWh]"
$ws.Range("M1").Value = "MAPE [%]"

# Update row 2 data values
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0
$ws.Range("D2").Value = 25
$ws.Range("G2").Value = "<keras.src.optimizers.adam.Adam object at 0x000001A05167D2A0>"
$ws.Range("I2").Value = 48
$ws.Range("J2").Value = 39.78723239730775
$ws.Range("K2").Value = 2230.414067202823
$ws.Range("L2").Value = 47.22725978926602
$ws.Range("M2").Value = 23.27835059092124
